$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cells I1 ("I0") and J1 ("IF") ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header formatting (bold font, border, centered alignment) from
# the existing H1 header cell onto the two new header cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Data columns I and J (rows 2-18) ---
$values = @(
    @(6, 7),
    @(6, 6),
    @(7, 7),
    @(7, 7),
    @(5, 6),
    @(5, 5),
    @(8, 8),
    @(9, 9),
    @(8, 8),
    @(5, 6),
    @(7, 8),
    @(10, 10),
    @(3, 4),
    @(6, 6),
    @(7, 7),
    @(5, 5),
    @(4, 4)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
